$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last data row dynamically (43 rows of data in this sheet)
$lastRow = $ws.UsedRange.Rows.Count

# Add header cells for the new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from the existing header cell (A1) to the new header cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in Wins/Losses/Ties values for every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD
    $ws.Cells.Item($r, 31).Value = 76   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
